$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BGI")

# Row 8
$ws.Range("D8").Value = 85100
$ws.Range("E8").Value = 86700
$ws.Range("F8").Value = 95700
$ws.Range("G8").Value = 106700
$ws.Range("H8").Value = 108900
$ws.Range("I8").Value = 217900
$ws.Range("J8").Value = 225000

# Row 9
$ws.Range("D9").Value = 52700
$ws.Range("E9").Value = 51800
$ws.Range("F9").Value = 56300
$ws.Range("G9").Value = 62700
$ws.Range("H9").Value = 61400
$ws.Range("I9").Value = 124000
$ws.Range("J9").Value = 125800

# Row 10
$ws.Range("D10").Value = 32400
$ws.Range("E10").Value = 34800
$ws.Range("F10").Value = 39400
$ws.Range("G10").Value = 44000
$ws.Range("H10").Value = 47400
$ws.Range("I10").Value = 93900
$ws.Range("J10").Value = 99200

# Row 14
$ws.Range("D14").Value = 2100
$ws.Range("E14").Value = 500
$ws.Range("F14").Value = 400
$ws.Range("G14").Value = 2700

# Row 15
$ws.Range("D15").Value = 1900
$ws.Range("E15").Value = 1900
$ws.Range("F15").Value = 2100
$ws.Range("G15").Value = 2300
$ws.Range("H15").Value = 2200
$ws.Range("I15").Value = 3400
$ws.Range("J15").Value = 3500

# Row 17
$ws.Range("D17").Value = 95300
$ws.Range("E17").Value = 89400
$ws.Range("F17").Value = 92400
$ws.Range("G17").Value = 107000
$ws.Range("H17").Value = 109000
$ws.Range("I17").Value = 209800
$ws.Range("J17").Value = 217200

# Row 18
$ws.Range("D18").Value = -10200
$ws.Range("E18").Value = -2800
$ws.Range("F18").Value = 3400
$ws.Range("G18").Value = -300
$ws.Range("H18").Value = -100
$ws.Range("I18").Value = 8000
$ws.Range("J18").Value = 7800

# Row 20
$ws.Range("D20").Value = 0

# Row 21
$ws.Range("D21").Value = -8300
$ws.Range("E21").Value = -800
$ws.Range("F21").Value = 5500
$ws.Range("G21").Value = 4100
$ws.Range("H21").Value = 4000
$ws.Range("I21").Value = 11600
$ws.Range("J21").Value = 11500

# Row 22
$ws.Range("D22").Value = 2300
$ws.Range("E22").Value = 2500
$ws.Range("F22").Value = 3200
$ws.Range("G22").Value = 2000
$ws.Range("H22").Value = 3700
$ws.Range("I22").Value = 6900
$ws.Range("J22").Value = 7600

# Row 23
$ws.Range("D23").Value = -12500
$ws.Range("E23").Value = -5300
$ws.Range("F23").Value = 200
$ws.Range("G23").Value = -2300
$ws.Range("H23").Value = -3800
$ws.Range("I23").Value = 1100

# Row 24
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 0

# Row 26
$ws.Range("D26").Value = -12500
$ws.Range("E26").Value = -5300
$ws.Range("F26").Value = 200
$ws.Range("G26").Value = -2300
$ws.Range("H26").Value = -3800
$ws.Range("I26").Value = 1100

# Row 27
$ws.Range("D27").Value = -12500
$ws.Range("E27").Value = -5300
$ws.Range("F27").Value = 200
$ws.Range("G27").Value = -2300
$ws.Range("H27").Value = -3800
$ws.Range("I27").Value = 1100

# Row 29
$ws.Range("D29").Value = 21200
$ws.Range("E29").Value = 8900
$ws.Range("F29").Value = 3900
$ws.Range("G29").Value = -4100
$ws.Range("H29").Value = -500
$ws.Range("I29").Value = "NA"
$ws.Range("J29").Value = "NA"

# Row 32
$ws.Range("D32").Value = 0

# Row 33
$ws.Range("D33").Value = 8700
$ws.Range("E33").Value = 3700
$ws.Range("F33").Value = 4000
$ws.Range("G33").Value = -6400
$ws.Range("H33").Value = -4300
$ws.Range("I33").Value = 1100

# Row 35
$ws.Range("D35").Value = 8700
$ws.Range("E35").Value = 3700
$ws.Range("F35").Value = 4000
$ws.Range("G35").Value = -6400
$ws.Range("H35").Value = -4300
$ws.Range("I35").Value = 1100

# Row 41
$ws.Range("D41").Value = 700
$ws.Range("E41").Value = 1400
$ws.Range("F41").Value = 1700
$ws.Range("G41").Value = 1800
$ws.Range("H41").Value = 1700
$ws.Range("I41").Value = 2800
$ws.Range("J41").Value = 2300

# Row 43
$ws.Range("D43").Value = 4600
$ws.Range("E43").Value = 10100
$ws.Range("F43").Value = 7700
$ws.Range("G43").Value = 5700
$ws.Range("H43").Value = 5500
$ws.Range("I43").Value = 4900
$ws.Range("J43").Value = 5300

# Row 44
$ws.Range("D44").Value = 63100
$ws.Range("E44").Value = 98300
$ws.Range("F44").Value = 102600
$ws.Range("G44").Value = 101000
$ws.Range("H44").Value = 107600
$ws.Range("I44").Value = 102000
$ws.Range("J44").Value = 106800

# Row 45
$ws.Range("D45").Value = 3700
$ws.Range("E45").Value = 1600
$ws.Range("F45").Value = 1300
$ws.Range("G45").Value = 1700
$ws.Range("H45").Value = 1800
$ws.Range("I45").Value = 1900
$ws.Range("J45").Value = 1800

# Row 46
$ws.Range("D46").Value = 72200
$ws.Range("E46").Value = 111500
$ws.Range("F46").Value = 113300
$ws.Range("G46").Value = 110200
$ws.Range("H46").Value = 116600
$ws.Range("I46").Value = 111700
$ws.Range("J46").Value = 116100

# Row 48
$ws.Range("D48").Value = 14500
$ws.Range("E48").Value = 17100
$ws.Range("F48").Value = 21900
$ws.Range("G48").Value = 21200
$ws.Range("H48").Value = 23000
$ws.Range("I48").Value = 20300
$ws.Range("J48").Value = 19200

# Row 49
$ws.Range("D49").Value = 2900
$ws.Range("E49").Value = 500
$ws.Range("F49").Value = 600
$ws.Range("G49").Value = 700
$ws.Range("H49").Value = 800
$ws.Range("I49").Value = 700
$ws.Range("J49").Value = 700

# Row 52
$ws.Range("D52").Value = 0
$ws.Range("E52").Value = 4100
$ws.Range("F52").Value = 400
$ws.Range("G52").Value = 2000
$ws.Range("H52").Value = 1400
$ws.Range("I52").Value = 1200
$ws.Range("J52").Value = 1700

# Row 54
$ws.Range("D54").Value = 89500
$ws.Range("E54").Value = 133200
$ws.Range("F54").Value = 136200
$ws.Range("G54").Value = 134100
$ws.Range("H54").Value = 141800
$ws.Range("I54").Value = 133900
$ws.Range("J54").Value = 137800

# Row 57
$ws.Range("D57").Value = 19600
$ws.Range("E57").Value = 34700
$ws.Range("F57").Value = 34800
$ws.Range("G57").Value = 33300
$ws.Range("H57").Value = 27500
$ws.Range("I57").Value = 31400
$ws.Range("J57").Value = 33200

# Row 58
$ws.Range("D58").Value = 30000
$ws.Range("E58").Value = 54500
$ws.Range("F58").Value = 50700
$ws.Range("G58").Value = 51400
$ws.Range("H58").Value = 58400
$ws.Range("I58").Value = 52900
$ws.Range("J58").Value = 50000

# Row 59
$ws.Range("D59").Value = 5600
$ws.Range("E59").Value = 6200
$ws.Range("F59").Value = 6700
$ws.Range("G59").Value = 6000
$ws.Range("H59").Value = 5900
$ws.Range("I59").Value = 6800
$ws.Range("J59").Value = 9700

# Row 60
$ws.Range("D60").Value = 55200
$ws.Range("E60").Value = 95500
$ws.Range("F60").Value = 92200
$ws.Range("G60").Value = 90700
$ws.Range("H60").Value = 91800
$ws.Range("I60").Value = 91200
$ws.Range("J60").Value = 92800

# Row 61
$ws.Range("D61").Value = 3600
$ws.Range("E61").Value = 22700
$ws.Range("F61").Value = 34700
$ws.Range("G61").Value = 38700
$ws.Range("H61").Value = 37400
$ws.Range("I61").Value = 28400
$ws.Range("J61").Value = 33800

# Row 62
$ws.Range("D62").Value = 6600
$ws.Range("E62").Value = 5500
$ws.Range("F62").Value = 3600
$ws.Range("G62").Value = 2600
$ws.Range("H62").Value = 2500
$ws.Range("I62").Value = 2200
$ws.Range("J62").Value = 2600

# Row 66
$ws.Range("D66").Value = 65300
$ws.Range("E66").Value = 123600
$ws.Range("F66").Value = 130400
$ws.Range("G66").Value = 132000
$ws.Range("H66").Value = 131600
$ws.Range("I66").Value = 121800
$ws.Range("J66").Value = 129100

# Row 72
$ws.Range("D72").Value = -59400
$ws.Range("E72").Value = -55000
$ws.Range("F72").Value = -58700
$ws.Range("G72").Value = -62700
$ws.Range("H72").Value = -56300
$ws.Range("I72").Value = -52000
$ws.Range("J72").Value = -53100

# Row 76
$ws.Range("D76").Value = 24200
$ws.Range("E76").Value = 9500
$ws.Range("F76").Value = 5700
$ws.Range("G76").Value = 2100
$ws.Range("H76").Value = 10100
$ws.Range("I76").Value = 12200
$ws.Range("J76").Value = 8700

# Row 81
$ws.Range("D81").Value = 8700
$ws.Range("E81").Value = 3700
$ws.Range("F81").Value = 4000
$ws.Range("G81").Value = -6400
$ws.Range("H81").Value = -4300
$ws.Range("I81").Value = 1100

# Row 83
$ws.Range("D83").Value = 1900
$ws.Range("E83").Value = 2000
$ws.Range("F83").Value = 2100
$ws.Range("G83").Value = 4500
$ws.Range("H83").Value = 4100
$ws.Range("I83").Value = 3600
$ws.Range("J83").Value = 3700

# Row 89
$ws.Range("D89").Value = -21800
$ws.Range("E89").Value = 5400
$ws.Range("F89").Value = 3500
$ws.Range("G89").Value = 7900
$ws.Range("H89").Value = -14200
$ws.Range("I89").Value = 4600
$ws.Range("J89").Value = 3400

# Row 91
$ws.Range("D91").Value = -3800
$ws.Range("E91").Value = -3300
$ws.Range("F91").Value = -3600
$ws.Range("G91").Value = -4700
$ws.Range("H91").Value = -4900
$ws.Range("I91").Value = -4700
$ws.Range("J91").Value = -3400

# Row 94
$ws.Range("D94").Value = 73600
$ws.Range("E94").Value = -3800
$ws.Range("F94").Value = -1800
$ws.Range("G94").Value = -4700
$ws.Range("H94").Value = -5100
$ws.Range("I94").Value = -4700
$ws.Range("J94").Value = -3400

# Row 100
$ws.Range("D100").Value = -52500
$ws.Range("E100").Value = -1900
$ws.Range("F100").Value = -1600
$ws.Range("G100").Value = -3000
$ws.Range("H100").Value = 18300
$ws.Range("I100").Value = 600
$ws.Range("J100").Value = -100

# Row 101
$ws.Range("D101").Value = -200
$ws.Range("F101").Value = -100
$ws.Range("G101").Value = -100
$ws.Range("H101").Value = -100

# Row 102
$ws.Range("D102").Value = -900
$ws.Range("E102").Value = -300
$ws.Range("H102").Value = -1100
$ws.Range("I102").Value = 500
